$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B11 value from 6 to 1
$ws.Range("B11").Value = 1

# Update the active cell selection to F11
$ws.Range("F11").Select()
